$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is unambiguous text (contains non-numeric
# characters such as a second "." or a "%" sign / spaces) can be
# written directly.
$textValues = @{
    "D2" = "26.644.28"
    "E2" = "  -0.19%  "
    "D3" = "1.596.59"
    "E3" = "  -0.05%  "
    "E4" = "  +0.20%  "
    "E5" = "  -0.13%  "
    "E6" = "  +0.56%  "
    "E7" = "  +0.18%  "
    "E8" = "  -0.12%  "
    "E9" = "  +0.43%  "
    "E10" = "  -1.20%  "
    "E11" = "  -0.25%  "
    "D12" = "1.821.18"
    "E12" = "  +0.00%  "
    "D13" = "1.654.64"
    "E13" = "  +3.62%  "
    "E14" = "  -0.15%  "
    "E15" = "  -0.04%  "
    "E16" = "  +0.07%  "
    "D17" = "26.622.72"
    "D18" = "0.0₃0737"
    "E18" = "  +1.15%  "
    "E19" = "  -0.18%  "
    "E21" = "  +3.53%  "
    "E22" = "  +0.40%  "
    "E23" = "  +1.33%  "
    "E24" = "  +0.66%  "
    "E25" = "  -1.20%  "
    "E26" = "  +0.10%  "
    "E27" = "  -1.13%  "
    "E29" = "  -0.31%  "
    "E30" = "  +2.09%  "
    "E31" = "  +0.28%  "
    "E32" = "  +0.76%  "
    "E33" = "  +1.28%  "
    "D34" = "1.281.80"
    "E34" = "  -1.16%  "
    "E35" = "  -7.09%  "
    "E36" = "  +0.62%  "
    "E37" = "  +0.87%  "
    "E38" = "  -0.90%  "
    "E39" = "  -1.21%  "
    "E40" = "  +19.84%  "
    "E41" = "  +2.28%  "
    "E42" = "  -0.06%  "
    "E43" = "  -0.75%  "
    "E44" = "  -0.10%  "
    "D45" = "1.734.29"
    "E45" = "  +0.04%  "
    "E46" = "  +0.49%  "
    "E47" = "  -3.21%  "
    "E48" = "  +1.60%  "
    "E49" = "  +0.66%  "
    "E50" = "  -0.11%  "
    "E51" = "  -1.72%  "
}

# Cells whose new value looks like a plain number (e.g. "19.53").
# Excel auto-converts such text to a real number on assignment, which
# would change the cell from a text/string cell to a numeric cell.
# The source sheet stores these as text, so force text interpretation
# (NumberFormat "@") for the assignment, then clear the format again so
# the cell keeps its original (default) style.
$numberLikeValues = @{
    "D5" = "211.34"
    "D10" = "19.53"
    "D16" = "65.12"
    "D19" = "209.83"
    "D23" = "2.34"
    "D25" = "144.95"
    "D31" = "1.16"
    "D35" = "0.620"
    "D41" = "5.50"
    "D44" = "63.72"
    "D46" = "90.59"
    "D49" = "0.0509"
    "D51" = "7.38"
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

foreach ($addr in $numberLikeValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numberLikeValues[$addr]
    $cell.ClearFormats()
}

